$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 was a blank placeholder row; fill it in with the new test case.
# Order of assignment matters for shared-string allocation order (B, then A, then C).
$ws.Range("B26").Value = "net is off, go to menu page, click on 'free bonus', ""could not connect to the server"" is showing, click on 'ok', app is crashing"
$ws.Range("A26").Value = "[windows-desktop]: app is crashing if clicking free bonus while net is off"
$ws.Range("C26").Value = "goto menu page after login, click on 'free bonus', click on 'ok'"
$ws.Range("D26").Value = 8306

# B26 formatting matches the style already used for long wrapped descriptions (e.g. B24).
$ws.Range("B26").VerticalAlignment = -4160
$ws.Range("B26").HorizontalAlignment = -4131
$ws.Range("B26").WrapText = $true
$ws.Range("B26").Borders.LineStyle = 1

# C26 gets a new left/top aligned (non-wrapped) style.
$ws.Range("C26").VerticalAlignment = -4160
$ws.Range("C26").HorizontalAlignment = -4131
$ws.Range("C26").Borders.LineStyle = 1

# Row grows to fit the new content.
$ws.Rows.Item(26).RowHeight = 93

# Leave the view positioned on the newly-added row/cell.
$ws.Range("C26").Select()
$excel.ActiveWindow.ScrollRow = 23
